# Swap the match-data (columns B:AC) between specific pairs of rows.
# Column A (the sequential row id) stays put; only the underlying
# match record (id, teams, odds, etc.) moves between the two rows in
# each pair, effectively re-ordering how same-date fixtures were listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Parallel arrays (PowerShell flattens arrays-of-arrays on assignment, so
# two flat arrays indexed in lock-step are used instead of an array of pairs).
$rowsA = @(50, 79, 98, 120, 193, 212, 220, 274, 296, 378, 394)
$rowsB = @(51, 80, 99, 121, 194, 213, 221, 275, 297, 379, 395)

for ($i = 0; $i -lt $rowsA.Count; $i++) {
    $r1 = $rowsA[$i]
    $r2 = $rowsB[$i]

    $range1 = $ws.Range("B" + $r1 + ":AC" + $r1)
    $range2 = $ws.Range("B" + $r2 + ":AC" + $r2)

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
